# Updates Leve profit-calculation figures across several Sheets, as produced
# by the scheduled market-data refresh runner. Each worksheet corresponds to
# a crafting discipline (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR); the edited
# rows are re-priced market cells (columns H-N).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ALC
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 709.14813
$ws.Range("J28").Value = 1098.2858
$ws.Range("L28").Value = 1098.2858
$ws.Range("N28").Value = -2068.2858

# ---------------------------------------------------------------------------
# BSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H86").Value = 56586.45
$ws.Range("I86").Value = 79897.86
$ws.Range("K86").Value = 79897.86
$ws.Range("M86").Value = -78774.86

$ws.Range("H89").Value = 56586.45
$ws.Range("I89").Value = 79897.86
$ws.Range("K89").Value = 399489.3
$ws.Range("M89").Value = -393873.3

$ws.Range("H135").Value = 39000
$ws.Range("J135").Value = 39000
$ws.Range("L135").Value = 39000
$ws.Range("N135").Value = -49140

# ---------------------------------------------------------------------------
# CRP
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H2").Value = 200
$ws.Range("I2").Value = 200
$ws.Range("K2").Value = 200
$ws.Range("M2").Value = -87

$ws.Range("H4").Value = 2083.3333
$ws.Range("I4").Value = 1500
$ws.Range("J4").Value = 5000
$ws.Range("K4").Value = 1500
$ws.Range("L4").Value = 5000
$ws.Range("M4").Value = -1388
$ws.Range("N4").Value = -5224

$ws.Range("H31").Value = 20655.207
$ws.Range("I31").Value = 24416.883
$ws.Range("J31").Value = 4480
$ws.Range("K31").Value = 24416.883
$ws.Range("L31").Value = 4480
$ws.Range("M31").Value = -24121.883
$ws.Range("N31").Value = -5070

$ws.Range("H34").Value = 20655.207
$ws.Range("I34").Value = 24416.883
$ws.Range("J34").Value = 4480
$ws.Range("K34").Value = 24416.883
$ws.Range("L34").Value = 4480
$ws.Range("M34").Value = -24214.883
$ws.Range("N34").Value = -4884

$ws.Range("H99").Value = 8936.5
$ws.Range("I99").Value = 2996.5
$ws.Range("J99").Value = 14876.5
$ws.Range("K99").Value = 2996.5
$ws.Range("L99").Value = 14876.5
$ws.Range("M99").Value = -1498.5
$ws.Range("N99").Value = -17872.5

$ws.Range("H126").Value = 8936.5
$ws.Range("I126").Value = 2996.5
$ws.Range("J126").Value = 14876.5
$ws.Range("K126").Value = 8989.5
$ws.Range("L126").Value = 44629.5
$ws.Range("M126").Value = -6519.5
$ws.Range("N126").Value = -49569.5

$ws.Range("H132").Value = 48389990
$ws.Range("I132").Value = 41669412
$ws.Range("J132").Value = 71431980
$ws.Range("K132").Value = 125008236
$ws.Range("L132").Value = 214295940
$ws.Range("M132").Value = -125005706
$ws.Range("N132").Value = -214301000

$ws.Range("H134").Value = 1931.8667
$ws.Range("I134").Value = 2456.25
$ws.Range("J134").Value = 1332.5714
$ws.Range("K134").Value = 7368.75
$ws.Range("L134").Value = 3997.7142
$ws.Range("M134").Value = -4833.75
$ws.Range("N134").Value = -9067.7142

# ---------------------------------------------------------------------------
# CUL
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H122").Value = 4131.778
$ws.Range("I122").Value = 330.22726
$ws.Range("K122").Value = 2972.04534
$ws.Range("M122").Value = -522.0453400000001

$ws.Range("H140").Value = 1299.4048
$ws.Range("I140").Value = 946.09375
$ws.Range("J140").Value = 2430
$ws.Range("K140").Value = 2838.28125
$ws.Range("L140").Value = 7290
$ws.Range("M140").Value = 2341.71875
$ws.Range("N140").Value = -17650

# ---------------------------------------------------------------------------
# GSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H70").Value = 88488.53999999999
$ws.Range("J70").Value = 5589.769
$ws.Range("L70").Value = 5589.769
$ws.Range("N70").Value = -6129.769

$ws.Range("H73").Value = 88488.53999999999
$ws.Range("J73").Value = 5589.769
$ws.Range("L73").Value = 5589.769
$ws.Range("N73").Value = -7461.769

$ws.Range("H126").Value = 2560106
$ws.Range("I126").Value = 2988
$ws.Range("K126").Value = 8964
$ws.Range("M126").Value = -6494

$ws.Range("H132").Value = 2147.2222
$ws.Range("I132").Value = 1526.8096
$ws.Range("K132").Value = 4580.4288
$ws.Range("M132").Value = -2050.4288

# ---------------------------------------------------------------------------
# LTW
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H2").Value = 236870.77
$ws.Range("I2").Value = 292233.6
$ws.Range("J2").Value = 104000
$ws.Range("K2").Value = 292233.6
$ws.Range("L2").Value = 104000
$ws.Range("M2").Value = -292121.6
$ws.Range("N2").Value = -104224

$ws.Range("H7").Value = 2211.3333
$ws.Range("I7").Value = 1740.8334
$ws.Range("J7").Value = 4093.3333
$ws.Range("K7").Value = 1740.8334
$ws.Range("L7").Value = 4093.3333
$ws.Range("M7").Value = -1628.8334
$ws.Range("N7").Value = -4317.3333

$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").Value = ""

$ws.Range("H40").Value = 68579.734
$ws.Range("J40").Value = 2226.9092
$ws.Range("L40").Value = 2226.9092
$ws.Range("N40").Value = -2498.9092

$ws.Range("H41").Value = 6934.25
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 6934.25
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 6934.25
$ws.Range("M41").Value = ""
$ws.Range("N41").Value = -7810.25

$ws.Range("H126").Value = 2211.3333
$ws.Range("I126").Value = 1740.8334
$ws.Range("J126").Value = 4093.3333
$ws.Range("K126").Value = 5222.5002
$ws.Range("L126").Value = 12279.9999
$ws.Range("M126").Value = -2752.5002
$ws.Range("N126").Value = -17219.9999

$ws.Range("H132").Value = 4055.818
$ws.Range("I132").Value = 5721.636
$ws.Range("J132").Value = 2390
$ws.Range("K132").Value = 17164.908
$ws.Range("L132").Value = 7170
$ws.Range("M132").Value = -14634.908
$ws.Range("N132").Value = -12230

# ---------------------------------------------------------------------------
# WVR
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H2").Value = 27907.143
$ws.Range("J2").Value = 34087.5
$ws.Range("L2").Value = 34087.5
$ws.Range("N2").Value = -34311.5

$ws.Range("H11").Value = 70005
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 70005
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 70005
$ws.Range("M11").Value = ""
$ws.Range("N11").Value = -70289

$ws.Range("H126").Value = 2212.9048
$ws.Range("I126").Value = 2342
$ws.Range("J126").Value = 1799.8
$ws.Range("K126").Value = 7026
$ws.Range("L126").Value = 5399.4
$ws.Range("M126").Value = -4556
$ws.Range("N126").Value = -10339.4

$ws.Range("H132").Value = 9056.857
$ws.Range("I132").Value = 15933.333
$ws.Range("K132").Value = 47799.999
$ws.Range("M132").Value = -45269.999
